$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.966.06"
$ws.Range("E2").Value = "  +0.94%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.639.64"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.00"
$ws.Range("E5").Value = "  +0.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.524"
$ws.Range("E6").Value = "  +0.26%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.56"
$ws.Range("E8").Value = "  +0.25%  "

# Row 9
$ws.Range("E9").Value = "  -1.31%  "

# Row 10
$ws.Range("E10").Value = "  +0.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("E11").Value = "  +2.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.67"
$ws.Range("E12").Value = "  +0.34%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.57"
$ws.Range("E13").Value = "  +0.05%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.575"
$ws.Range("E14").Value = "  +3.79%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.10"
$ws.Range("E15").Value = "  +1.11%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.91"
$ws.Range("E16").Value = "  +1.17%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.968.99"
$ws.Range("E17").Value = "  +1.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.01"
$ws.Range("E18").Value = "  +1.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +0.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("E20").Value = "  +0.12%  "

# Row 21
$ws.Range("E21").Value = "  +0.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("E22").Value = "  +0.68%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.36"
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
$ws.Range("E24").Value = "  -1.49%  "

# Row 25
$ws.Range("E25").Value = "  +1.47%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.00"
$ws.Range("E26").Value = "  +1.42%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.70"
$ws.Range("E27").Value = "  +0.77%  "

# Row 29
$ws.Range("E29").Value = "  +0.17%  "

# Row 30
$ws.Range("E30").Value = "  +0.18%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0483"
$ws.Range("E31").Value = "  -0.14%  "

# Row 32
$ws.Range("E32").Value = "  +1.76%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.12"
$ws.Range("E33").Value = "  +1.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.406.99"
$ws.Range("E34").Value = "  -5.16%  "

# Row 35
$ws.Range("E35").Value = "  +2.09%  "

# Row 36
$ws.Range("E36").Value = "  +1.08%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0169"
$ws.Range("E37").Value = "  +0.89%  "

# Row 38
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.882"
$ws.Range("E38").Value = "  -0.09%  "

# Row 39
$ws.Range("E39").Value = "  -0.06%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.906"
$ws.Range("E40").Value = "  -5.99%  "

# Row 41
$ws.Range("E41").Value = "  +0.41%  "

# Row 42
$ws.Range("E42").Value = "  +0.06%  "

# Row 43
$ws.Range("E43").Value = "  +6.96%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.52"
$ws.Range("E44").Value = "  +3.50%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.33"
$ws.Range("E45").Value = "  -2.24%  "

# Row 46
$ws.Range("E46").Value = "  -0.10%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.781.27"
$ws.Range("E47").Value = "  +0.43%  "

# Row 48
$ws.Range("E48").Value = "  +0.44%  "

# Row 49
$ws.Range("E49").Value = "  +1.44%  "

# Row 50
$ws.Range("E50").Value = "  +0.31%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  -1.07%  "
